# Sort the data rows (A5:H8) ascending by column A (Cylinder), then
# merge the two rows that end up sharing the same Cylinder value (6)
# into a single spanning cell, top-aligned, leaving the lower half of
# the merge blank (matching the formatting of the other blank "filler"
# cells in column H).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sort rows 5-8 (the data body) by column A, ascending.
$dataRange = $ws.Range("A5:H8")
$sortKey = $ws.Range("A5:A8")
$dataRange.Sort($sortKey)

# After sorting, rows 6 and 7 both have Cylinder == 6; merge A6:A7.
$ws.Range("A6:A7").Merge()

# The merged (visible) cell keeps the bold "value" style, vertically
# top-aligned.
$ws.Range("A6").VerticalAlignment = -4160

# The now-blank lower half of the merge (A7) should look like the
# other blank filler cells (e.g. H5) rather than keep the bold style.
$ws.Range("H5").Copy()
$ws.Range("A7").PasteSpecial(-4122)
